$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.381.28'
$ws.Range('E2').Value = '  +5.52%  '

$ws.Range('D3').Value = '3.281.11'
$ws.Range('E3').Value = '  +2.78%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.13'
$ws.Range('E5').Value = '  -2.89%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '629.91'
$ws.Range('E6').Value = '  +1.56%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.386'
$ws.Range('E7').Value = '  +29.67%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.671'
$ws.Range('E8').Value = '  +14.73%  '

$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').Value = '3.275.91'
$ws.Range('E10').Value = '  +2.96%  '

$ws.Range('E11').Value = '  -2.22%  '

$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000261'
$ws.Range('E12').Value = '  -0.82%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.178'
$ws.Range('E13').Value = '  +7.71%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.34'
$ws.Range('E14').Value = '  +6.43%  '

$ws.Range('D15').Value = '3.882.47'
$ws.Range('E15').Value = '  +3.04%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.33'
$ws.Range('E16').Value = '  +0.17%  '

$ws.Range('D17').Value = '87.110.34'
$ws.Range('E17').Value = '  +5.57%  '

$ws.Range('D18').Value = '3.273.59'
$ws.Range('E18').Value = '  +3.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.12'
$ws.Range('E19').Value = '  +0.17%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.05'
$ws.Range('E20').Value = '  -7.56%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.34'
$ws.Range('E21').Value = '  -1.57%  '

$ws.Range('E22').Value = '  -0.31%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.32'
$ws.Range('E23').Value = '  +3.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.16'
$ws.Range('E24').Value = '  -2.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.49'
$ws.Range('E25').Value = '  +4.27%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.16'
$ws.Range('E26').Value = '  -1.49%  '

$ws.Range('D27').Value = '3.448.41'
$ws.Range('E27').Value = '  +3.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '76.32'
$ws.Range('E28').Value = '  -2.14%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000131'
$ws.Range('E29').Value = '  +7.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.178'
$ws.Range('E31').Value = '  +18.93%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.61%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.83'
$ws.Range('E33').Value = '  -3.22%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '548.15'
$ws.Range('E34').Value = '  -3.63%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.43'
$ws.Range('E35').Value = '  -3.98%  '

$ws.Range('E36').Value = '  -1.83%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.03'
$ws.Range('E37').Value = '  +12.71%  '

$ws.Range('E38').Value = '  -10.06%  '

$ws.Range('E39').Value = '  -0.96%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.07%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '21.63'
$ws.Range('E41').Value = '  +3.71%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.396'
$ws.Range('E42').Value = '  -2.42%  '

$ws.Range('E43').Value = '  -1.01%  '

$ws.Range('E44').Value = '  -1.65%  '

$ws.Range('E45').Value = '  +0.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '156.68'
$ws.Range('E46').Value = '  -2.52%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '179.28'
$ws.Range('E47').Value = '  -4.06%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '44.62'
$ws.Range('E48').Value = '  -0.14%  '

$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('E50').Value = '  +1.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.627'
$ws.Range('E51').Value = '  -0.57%  '
